# Fruta / hortaliza, semanal
# Insert two new weekly price records at the top of the data block
# (rows 230:231), pushing the existing rows 230:247 down to 232:249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 230 (Excel shifts
# everything at/after row 230 down by 2, i.e. old 230..247 -> 232..249).
$ws.Rows("230:231").Insert()

# New record 1 (row 230): Zafiro rojo
$ws.Range("A230").Value = 11
$ws.Range("B230").Value = "Vega Monumental Concepción"
$ws.Range("C230").Value = "Bíobío"
$ws.Range("D230").Value = 44610
$ws.Range("E230").Value = 8
$ws.Range("F230").Value = 100112002
$ws.Range("G230").Value = "Pimiento"
$ws.Range("H230").Value = "Zafiro rojo"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 100
$ws.Range("K230").Value = 23000
$ws.Range("L230").Value = 25000
$ws.Range("M230").Value = 24000
$ws.Range("N230").Value = "$/caja 15 kilos"
$ws.Range("O230").Value = "Región de Arica y Parinacota"
$ws.Range("P230").Value = 1600
$ws.Range("Q230").Value = 15
$ws.Range("R230").Value = "Hortaliza"

# New record 2 (row 231): Zafiro verde
$ws.Range("A231").Value = 11
$ws.Range("B231").Value = "Vega Monumental Concepción"
$ws.Range("C231").Value = "Bíobío"
$ws.Range("D231").Value = 44610
$ws.Range("E231").Value = 8
$ws.Range("F231").Value = 100112002
$ws.Range("G231").Value = "Pimiento"
$ws.Range("H231").Value = "Zafiro verde"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 100
$ws.Range("K231").Value = 13000
$ws.Range("L231").Value = 15000
$ws.Range("M231").Value = 14000
$ws.Range("N231").Value = "$/caja 15 kilos"
$ws.Range("O231").Value = "Región de Arica y Parinacota"
$ws.Range("P231").Value = 933
$ws.Range("Q231").Value = 15
$ws.Range("R231").Value = "Hortaliza"
